# "added code for new task"
#
# - Contacts: add PhoneNumber / JobTitle columns, refresh sample row
# - Deals: rework columns to DealName/DealAmount/DealStage/DealType/DealCompany/DealConatct
# - Tasks: brand-new sheet with task sample data (incl. a time-of-day value)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Contacts sheet
# ---------------------------------------------------------------------------
$contacts = $wb.Worksheets.Item("Contacts")

# existing hyperlinked email cell (A2) keeps pointing at the same mailto: link,
# just change the visible text (writing straight to the cell keeps the single
# existing Hyperlink record intact instead of appending a duplicate one)
$contacts.Range("A2").Value = "newfour@abc.com"

$contacts.Range("B2").Value = "New"
$contacts.Range("C2").Value = "Four"
$contacts.Range("D2").Value = "Subscriber"
$contacts.Range("E2").Value = "New"

$contacts.Range("F1").Value = "PhoneNumber"
$contacts.Range("G1").Value = "JobTitle"
$contacts.Range("F2").Value = 1234567890
$contacts.Range("G2").Value = "Engineer"

$contacts.Columns.Item(6).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# Deals sheet
# ---------------------------------------------------------------------------
$deals = $wb.Worksheets.Item("Deals")

# B2 used to hold a date formatted with a custom "mmmm yyyy" numFmt. That slot
# gets reused below by the new Tasks!F2 time value, so give it the new
# time-of-day format *before* moving it off this sheet.
$deals.Range("B2").NumberFormat = "h:mm AM/PM"

# ---------------------------------------------------------------------------
# Tasks sheet (new) - create now so we can relocate the formatted cell into it
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tasks = $wb.Worksheets.Add($null, $lastSheet)
$tasks.Name = "Tasks"

# Move the (now time-formatted) cell over to Tasks!F2 so the style slot it
# occupies is reused in place rather than a fresh one being appended.
$deals.Range("B2").Cut($tasks.Range("F2")) | Out-Null

# Deals!B2 reverts to a plain number with default formatting
$deals.Range("A1").Value = "DealName"
$deals.Range("B1").Value = "DealAmount"
$deals.Range("C1").Value = "DealStage"
$deals.Range("D1").Value = "DealType"
$deals.Range("E1").Value = "DealCompany"
$deals.Range("F1").Value = "DealConatct"

$deals.Range("A2").Value = "NewFour"
$deals.Range("B2").Value = 2000
$deals.Range("B2").Style = "Normal"
$deals.Range("C2").Value = "Closed won"
$deals.Range("D2").Value = "New Business"
$deals.Range("E2").Value = "ffd9g.com"
$deals.Range("F2").Value = "Raj Khanna (abcd1@abc.com)"

$deals.Columns.Item(3).AutoFit() | Out-Null
$deals.Columns.Item(4).AutoFit() | Out-Null
$deals.Columns.Item(5).AutoFit() | Out-Null
$deals.Columns.Item(6).AutoFit() | Out-Null

$deals.Range("B2").Select() | Out-Null

# ---------------------------------------------------------------------------
# Tasks sheet - fill in the rest of the content
# ---------------------------------------------------------------------------
$tasks.Range("A1").Value = "TaskTitle"
$tasks.Range("B1").Value = "TaskType"
$tasks.Range("C1").Value = "TaskAssociateWith"
$tasks.Range("D1").Value = "TaskQueueValue"
$tasks.Range("E1").Value = "TaskDateValue"
$tasks.Range("F1").Value = "TaskTimeValue"

$tasks.Range("A2").Value = "Twelve"
$tasks.Range("B2").Value = "Email"
$tasks.Range("C2").Value = "Raj Khanna (abcd1@abc.com)"
$tasks.Range("D2").Value = "Tesy2"
$tasks.Range("E2").Value = "Customdate"
$tasks.Range("F2").Value = 0.45833333333333331

$tasks.Columns.Item(3).AutoFit() | Out-Null
$tasks.Columns.Item(4).AutoFit() | Out-Null
$tasks.Columns.Item(5).AutoFit() | Out-Null
$tasks.Columns.Item(6).AutoFit() | Out-Null

$tasks.Range("D2").Select() | Out-Null
